$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.417.07'
$ws.Range("E2").Value = '  +1.12%  '

$ws.Range("D3").Value = '1.666.78'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.46'
$ws.Range("E5").Value = '  +1.72%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3968'
$ws.Range("E7").Value = '  +1.50%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3921'
$ws.Range("E8").Value = '  +1.51%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.08'
$ws.Range("E9").Value = '  +6.92%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.403'
$ws.Range("E10").Value = '  +3.16%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.001'
$ws.Range("E11").Value = '  -0.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08621'
$ws.Range("E12").Value = '  +1.81%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.43'
$ws.Range("E13").Value = '  +1.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.345'
$ws.Range("E14").Value = '  +2.65%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001360'
$ws.Range("E15").Value = '  +5.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.922'
$ws.Range("E16").Value = '  +5.71%  '

$ws.Range("D17").Value = '1.663.86'
$ws.Range("E17").Value = '  +1.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.48'
$ws.Range("E18").Value = '  +1.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06988'
$ws.Range("E19").Value = '  +0.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.62'
$ws.Range("E20").Value = '  -1.72%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.020'
$ws.Range("E21").Value = '  +0.90%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9994'
$ws.Range("E22").Value = '  -0.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.79'
$ws.Range("E23").Value = '  +0.40%  '

$ws.Range("D24").Value = '24.427.03'
$ws.Range("E24").Value = '  +1.15%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.427'
$ws.Range("E25").Value = '  +3.40%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.040'
$ws.Range("E26").Value = '  +11.37%  '

$ws.Range("E27").Value = '  +0.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.76'

$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '142.87'
$ws.Range("E29").Value = '  +0.78%  '

$ws.Range("B30").Value = 'HuobiToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.454'
$ws.Range("E30").Value = '  +1.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.124'
$ws.Range("E31").Value = '  -9.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.509'
$ws.Range("E32").Value = '  +1.44%  '

$ws.Range("D33").Value = '1.845.60'
$ws.Range("E33").Value = '  +1.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.075'
$ws.Range("E34").Value = '  +9.35%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08302'
$ws.Range("E35").Value = '  +3.40%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.03037'
$ws.Range("E36").Value = '  +3.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.918'
$ws.Range("E37").Value = '  -4.42%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2781'
$ws.Range("E38").Value = '  +2.56%  '

$ws.Range("B39").Value = 'FraxShare'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.14'
$ws.Range("E39").Value = '  +10.69%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09259'
$ws.Range("E40").Value = '  +0.07%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7754'
$ws.Range("E41").Value = '  +1.51%  '

$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.88'
$ws.Range("E42").Value = '  +5.67%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.445'
$ws.Range("E43").Value = '  -2.24%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.66'
$ws.Range("E44").Value = '  +4.27%  '

$ws.Range("E45").Value = '  +3.58%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.548'
$ws.Range("E46").Value = '  +2.18%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.145'
$ws.Range("E47").Value = '  +1.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  -0.28%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08464'
$ws.Range("E49").Value = '  +0.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.74'
$ws.Range("E50").Value = '  +1.90%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.280'
$ws.Range("E51").Value = '  +1.20%  '
